# "this is second mod"
# Adds a new row (row 3) to Sheet1 that duplicates row 2's layout/styling
# but with the "UE" label in column B instead of "DL", then selects B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (values + formatting) into row 3.
$ws.Range("A2:E2").Copy($ws.Range("A3:E3"))

# Column B of the new row should read "UE" (new shared string).
$ws.Range("B3").Value = "UE"

# Match the saved selection state from the edit.
$ws.Range("B3").Select() | Out-Null
